# Incorporate updated data from upstream processes through 2024
#
# The "Solar" facilities-opened counts for 2022 (row 24) and 2024
# (row 26) change from 35 -> 36 and 43 -> 64 respectively. The
# embedded column chart's Solar series reads its values directly from
# this same worksheet range (Sheet1!$E$2:$E$26), so updating the two
# source cells is the authoritative edit - the chart reflects the new
# numbers from these cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E24").Value = 36
$ws.Range("E26").Value = 64
